$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("B2").Value = 34.766589115695034
$ws.Range("C2").Value = 15.964111878738303
$ws.Range("D2").Value = 0.4591796976578143
$ws.Range("E2").Value = 30.772227766710479
$ws.Range("F2").Value = 15.367492730672387
$ws.Range("G2").Value = 0.49939487147878853
$ws.Range("H2").Value = 338.5
$ws.Range("I2").Value = 302.5
$ws.Range("B3").Value = 34.811697109202491
$ws.Range("C3").Value = 16.005152026127682
$ws.Range("D3").Value = 0.45976362416116484
$ws.Range("E3").Value = 30.812307935858165
$ws.Range("F3").Value = 15.405672362614842
$ws.Range("G3").Value = 0.49998436970981719
$ws.Range("H3").Value = 338.5
$ws.Range("I3").Value = 302.5
$ws.Range("B4").Value = 34.838277613680795
$ws.Range("C4").Value = 16.03046497794686
$ws.Range("D4").Value = 0.4601394235302777
$ws.Range("E4").Value = 30.850187059445034
$ws.Range("F4").Value = 15.449652511215708
$ws.Range("G4").Value = 0.50079607236889256
$ws.Range("H4").Value = 338.5
$ws.Range("I4").Value = 302
$ws.Range("B5").Value = 34.868685727598681
$ws.Range("C5").Value = 16.056062460222176
$ws.Range("D5").Value = 0.46047225827940363
$ws.Range("E5").Value = 30.880684803410553
$ws.Range("F5").Value = 15.477324874711522
$ws.Range("G5").Value = 0.50119759238636319
$ws.Range("H5").Value = 338.5
$ws.Range("I5").Value = 302
$ws.Range("B6").Value = 34.877566029482153
$ws.Range("C6").Value = 16.064873538813725
$ws.Range("D6").Value = 0.46060764461700165
$ws.Range("E6").Value = 30.908911721550538
$ws.Range("F6").Value = 15.49151019126575
$ws.Range("G6").Value = 0.50119882352456446
$ws.Range("H6").Value = 338.5
$ws.Range("I6").Value = 302
$ws.Range("B7").Value = 34.889874079538188
$ws.Range("C7").Value = 16.073750341553687
$ws.Range("D7").Value = 0.46069958019654861
$ws.Range("E7").Value = 30.933795506608398
$ws.Range("F7").Value = 15.51478705612951
$ws.Range("G7").Value = 0.50154812243505908
$ws.Range("H7").Value = 338.5
$ws.Range("I7").Value = 302
$ws.Range("B8").Value = 34.899093733979356
$ws.Range("C8").Value = 16.082601439180841
$ws.Range("D8").Value = 0.46083149212330649
$ws.Range("E8").Value = 30.950956582531738
$ws.Range("F8").Value = 15.525912020710509
$ws.Range("G8").Value = 0.50162947239805522
$ws.Range("H8").Value = 338.5
$ws.Range("I8").Value = 302
$ws.Range("B9").Value = 34.91417057195526
$ws.Range("C9").Value = 16.09713715254377
$ws.Range("D9").Value = 0.46104882025963878
$ws.Range("E9").Value = 30.639697009053023
$ws.Range("F9").Value = 15.262229435687013
$ws.Range("G9").Value = 0.49811946349135
$ws.Range("H9").Value = 338.5
$ws.Range("I9").Value = 301
$ws.Range("B10").Value = 34.925126585070657
$ws.Range("C10").Value = 16.107087224092599
$ws.Range("D10").Value = 0.46118908645496071
$ws.Range("E10").Value = 30.651034663778006
$ws.Range("F10").Value = 15.276738165447519
$ws.Range("G10").Value = 0.49840856378988313
$ws.Range("H10").Value = 338.5
$ws.Range("I10").Value = 301
$ws.Range("B11").Value = 34.938369452567947
$ws.Range("C11").Value = 16.116278001432512
$ws.Range("D11").Value = 0.4612773364627632
$ws.Range("E11").Value = 31.017645893353595
$ws.Range("F11").Value = 15.616666433191632
$ws.Range("G11").Value = 0.50347684304881246
$ws.Range("H11").Value = 338.5
$ws.Range("I11").Value = 301.5
